$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as text (inline strings) in the
# source workbook. Using a leading apostrophe forces Excel to keep the
# assigned value as text instead of auto-converting it to a number,
# exactly like typing '246.29 into a cell while it mirrors the original
# data update produced by the price-scraper GitHub Action.
$ws.Range("D2").Value = "'246.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.317"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05878"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.374"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8129"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9568"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.03703"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.07340"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03055"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'4.417"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09398"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001605"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04805"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005902"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "'0.004083"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0009894"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00009709"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Value = "'0.1284"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002472"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03892"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006756"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'0.002442"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.005912"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005670"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.6517"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.05609"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
